$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving original (unstyled) appearance
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 51: coin renamed from EnergySwap to dogwifhat, with new link, price and volume
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "3.30"
Set-TextValue "E51" "  +7.97%  "

# Price (column D) updates
Set-TextValue "D2" "97.103.46"
Set-TextValue "D3" "3.572.23"
Set-TextValue "D5" "241.19"
Set-TextValue "D6" "654.57"
Set-TextValue "D11" "3.568.60"
Set-TextValue "D12" "44.16"
Set-TextValue "D13" "0.204"
Set-TextValue "D15" "4.236.32"
Set-TextValue "D16" "96.856.28"
Set-TextValue "D17" "0.0000260"
Set-TextValue "D19" "3.571.61"
Set-TextValue "D20" "12.70"
Set-TextValue "D21" "17.97"
Set-TextValue "D22" "0.528"
Set-TextValue "D23" "3.48"
Set-TextValue "D24" "513.37"
Set-TextValue "D26" "6.85"
Set-TextValue "D27" "101.59"
Set-TextValue "D28" "13.05"
Set-TextValue "D29" "3.763.68"
Set-TextValue "D30" "0.168"
Set-TextValue "D32" "11.89"
Set-TextValue "D34" "0.184"
Set-TextValue "D35" "0.998"
Set-TextValue "D36" "31.76"
Set-TextValue "D37" "8.86"
Set-TextValue "D38" "616.71"
Set-TextValue "D40" "1.65"
Set-TextValue "D41" "1.96"
Set-TextValue "D44" "0.922"
Set-TextValue "D45" "6.02"
Set-TextValue "D46" "0.0438"
Set-TextValue "D47" "2.30"
Set-TextValue "D48" "23.60"

# Volume(1h) (column E) updates
Set-TextValue "E3" "  +0.06%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "E5" "  +2.50%  "
Set-TextValue "E6" "  +0.21%  "
Set-TextValue "E7" "  +16.63%  "
Set-TextValue "E8" "  +6.44%  "
Set-TextValue "E9" "  -0.07%  "
Set-TextValue "E10" "  +4.73%  "
Set-TextValue "E11" "  +0.04%  "
Set-TextValue "E12" "  +4.23%  "
Set-TextValue "E13" "  +0.41%  "
Set-TextValue "E14" "  -0.95%  "
Set-TextValue "E15" "  -0.05%  "
Set-TextValue "E16" "  +1.94%  "
Set-TextValue "E17" "  +2.64%  "
Set-TextValue "E18" "  +11.41%  "
Set-TextValue "E19" "  -0.28%  "
Set-TextValue "E20" "  +0.79%  "
Set-TextValue "E21" "  +0.79%  "
Set-TextValue "E22" "  +10.20%  "
Set-TextValue "E23" "  +0.91%  "
Set-TextValue "E24" "  +0.94%  "
Set-TextValue "E25" "  +5.63%  "
Set-TextValue "E26" "  +1.24%  "
Set-TextValue "E27" "  +6.68%  "
Set-TextValue "E28" "  +2.66%  "
Set-TextValue "E29" "  +0.04%  "
Set-TextValue "E30" "  +17.61%  "
Set-TextValue "E31" "  -1.32%  "
Set-TextValue "E32" "  +3.20%  "
Set-TextValue "E33" "  -0.14%  "
Set-TextValue "E34" "  +3.43%  "
Set-TextValue "E35" "  +0.04%  "
Set-TextValue "E36" "  +0.02%  "
Set-TextValue "E37" "  +4.23%  "
Set-TextValue "E38" "  +5.86%  "
Set-TextValue "E39" "  +1.22%  "
Set-TextValue "E40" "  -2.18%  "
Set-TextValue "E41" "  +7.83%  "
Set-TextValue "E42" "  +2.42%  "
Set-TextValue "E43" "  -0.10%  "
Set-TextValue "E44" "  +1.82%  "
Set-TextValue "E45" "  +4.74%  "
Set-TextValue "E46" "  +5.61%  "
Set-TextValue "E47" "  +0.52%  "
Set-TextValue "E48" "  +0.91%  "
Set-TextValue "E49" "  +31.52%  "
Set-TextValue "E50" "  +3.94%  "
